$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 363) holds a date serial value that was
# bumped forward by exactly one day (46075 -> 46076) during the
# automatic update.
$ws.Range("C2:C363").Value = 46076
